$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update "Datos actualizados" timestamp in A1 (06:05 -> 07:05)
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 07:05"

# 2. Swap rows 200 <-> 201 (Belice / Santa Lucia reorder), columns A:H
for ($c = 1; $c -le 8; $c++) {
    $v200 = $ws.Cells.Item(200, $c).Value2
    $v201 = $ws.Cells.Item(201, $c).Value2
    $ws.Cells.Item(200, $c).Value = $v201
    $ws.Cells.Item(201, $c).Value = $v200
}

# 3. Swap rows 213 <-> 214 (Islas Virgenes Britanicas / Papua Nueva Guinea reorder), columns A:H
for ($c = 1; $c -le 8; $c++) {
    $v213 = $ws.Cells.Item(213, $c).Value2
    $v214 = $ws.Cells.Item(214, $c).Value2
    $ws.Cells.Item(213, $c).Value = $v214
    $ws.Cells.Item(214, $c).Value = $v213
}

# 4. Update India (row 12) figures
$ws.Range("B12").Value = 182143
$ws.Range("C12").Value = 316
$ws.Range("D12").Value = 86984
$ws.Range("E12").Value = 89974

# 5. Update Tailandia (row 81) figures
$ws.Range("B81").Value = 3081
$ws.Range("C81").Value = 4
$ws.Range("D81").Value = 2963
$ws.Range("E81").Value = 61

# 6. Update Kirguistan (row 97) figures
$ws.Range("B97").Value = 1748
$ws.Range("C97").Value = 26
$ws.Range("D97").Value = 1170
$ws.Range("E97").Value = 562
